# Applies updated betting odds to row 2 and row 8 of Sheet1,
# matching the "Atualizando o arquivo XLSX" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value  = 1.4
$ws.Range("H2").Value  = 4.5
$ws.Range("I2").Value  = 7.5
$ws.Range("J2").Value  = 1.95
$ws.Range("K2").Value  = 2.25
$ws.Range("L2").Value  = 8
$ws.Range("M2").Value  = 1.05
$ws.Range("N2").Value  = 11
$ws.Range("O2").Value  = 1.3
$ws.Range("P2").Value  = 3.4
$ws.Range("Q2").Value  = 2
$ws.Range("R2").Value  = 1.8
$ws.Range("U2").Value  = 2.25
$ws.Range("V2").Value  = 1.57
$ws.Range("W2").Value  = 5.5
$ws.Range("X2").Value  = 6
$ws.Range("Z2").Value  = 8.5
$ws.Range("AD2").Value = 9
$ws.Range("AE2").Value = 23
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 41
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 101
$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 67
$ws.Range("AN2").Value = 3.2
$ws.Range("AO2").Value = 7
$ws.Range("AQ2").Value = 21
$ws.Range("AS2").Value = 201
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 8.5
$ws.Range("AZ2").Value = 201

# --- Row 8 updates ---
$ws.Range("R8").Value  = 1.57
$ws.Range("BD8").Value = 151
